$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.222.52"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.246.18"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'243.44"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'74.40"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").Value = "'42.31"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").Value = "'0.0958"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "'6.96"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "2.582.52"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "2.239.87"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "42.150.59"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "'72.83"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Value = "'11.18"
$ws.Range("E22").Value = "  +8.81%  "
$ws.Range("D23").Value = "'230.57"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("E24").Value = "  -6.13%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "'167.90"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'20.63"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("E32").Value = "  -4.49%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "'30.00"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -5.30%  "
$ws.Range("D37").Value = "'4.33"
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("D38").Value = "'0.0305"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").Value = "'13.19"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "'104.51"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").Value = "2.455.93"
$ws.Range("E51").Value = "  +0.36%  "

# Restore default (unstyled) formatting on cells that required a
# text-forcing apostrophe prefix, so no stray quotePrefix style sticks.
$ws.Range("D5").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D45").ClearFormats()
